$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up / normalize phone number formatting in column D
$ws.Range("D8").Value = "(416) 508-1780"
$ws.Range("D13").Value = "(778) 709-1769"
$ws.Range("D15").Value = "(778) 846-9554"
$ws.Range("D16").Value = "(250) 386-5311-2030"
$ws.Range("D17").Value = "(250) 386-5311-2030"
